$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.811.83"
$ws.Range("E2").Value = "  -4.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.458.70"
$ws.Range("E3").Value = "  -5.91%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "546.73"
$ws.Range("E5").Value = "  -4.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.00"
$ws.Range("E6").Value = "  -7.19%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.594"
$ws.Range("E8").Value = "  -4.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.456.68"
$ws.Range("E9").Value = "  -5.91%  "
$ws.Range("E10").Value = "  -9.80%  "
$ws.Range("E11").Value = "  -1.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.33"
$ws.Range("E12").Value = "  -8.62%  "
$ws.Range("E13").Value = "  -8.15%  "
$ws.Range("E14").Value = "  -8.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.899.24"
$ws.Range("E15").Value = "  -6.03%  "
$ws.Range("E16").Value = "  -10.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.733.73"
$ws.Range("E17").Value = "  -4.60%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.457.37"
$ws.Range("E18").Value = "  -6.07%  "
$ws.Range("E19").Value = "  -8.51%  "
$ws.Range("E20").Value = "  -8.30%  "
$ws.Range("E21").Value = "  -8.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "317.45"
$ws.Range("E22").Value = "  -7.56%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.03"
$ws.Range("E24").Value = "  -6.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.73"
$ws.Range("E25").Value = "  -5.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0₃0980"
$ws.Range("E26").Value = "  -9.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.576.46"
$ws.Range("E27").Value = "  -5.67%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  -5.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "528.68"
$ws.Range("E30").Value = "  -11.66%  "
$ws.Range("E31").Value = "  -9.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.60"
$ws.Range("E32").Value = "  -4.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.148"
$ws.Range("E33").Value = "  -8.34%  "
$ws.Range("E34").Value = "  -8.78%  "
$ws.Range("E35").Value = "  -9.79%  "
$ws.Range("E36").Value = "  -11.86%  "
$ws.Range("E37").Value = "  -10.42%  "
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("E39").Value = "  -6.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.26"
$ws.Range("E40").Value = "  -7.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "144.55"
$ws.Range("E41").Value = "  -6.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.70"
$ws.Range("E43").Value = "  -9.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.90"
$ws.Range("E44").Value = "  -4.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.27"
$ws.Range("E45").Value = "  -9.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "146.47"
$ws.Range("E46").Value = "  -6.85%  "
$ws.Range("E47").Value = "  -8.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.77"
$ws.Range("E48").Value = "  -12.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0527"
$ws.Range("E49").Value = "  -10.85%  "

# Row 50/51 swap: Mantle <-> Stellar
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0941"
$ws.Range("E50").Value = "  -6.13%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.580"
$ws.Range("E51").Value = "  -7.98%  "
